# Updates cryptos list prices/volumes (and Monero/Toncoin row swap) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    # Prefix with an apostrophe so Excel stores numeric-looking strings
    # (e.g. "0.519", "2.20") as text instead of converting them to numbers,
    # then reset the style so the quote-prefix flag does not linger.
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).Style = "Normal"
}

Set-TextCell "D2" "43.019.87"
Set-TextCell "E2" "  +1.83%  "
Set-TextCell "D3" "2.308.23"
Set-TextCell "E3" "  +1.55%  "
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.08%  "
Set-TextCell "D5" "302.88"
Set-TextCell "E5" "  +1.01%  "
Set-TextCell "D6" "101.25"
Set-TextCell "E6" "  +5.10%  "
Set-TextCell "E7" "  +1.69%  "
Set-TextCell "E8" "  -0.12%  "
Set-TextCell "D9" "0.519"
Set-TextCell "E9" "  +5.18%  "
Set-TextCell "D10" "35.50"
Set-TextCell "E10" "  +6.50%  "
Set-TextCell "D11" "0.0796"
Set-TextCell "E11" "  +1.00%  "
Set-TextCell "E12" "  +3.79%  "
Set-TextCell "D13" "18.05"
Set-TextCell "E13" "  +15.76%  "
Set-TextCell "D14" "6.93"
Set-TextCell "E14" "  +3.78%  "
Set-TextCell "D15" "2.684.43"
Set-TextCell "E15" "  +2.16%  "
Set-TextCell "D16" "2.314.39"
Set-TextCell "E16" "  +1.07%  "
Set-TextCell "D17" "0.814"
Set-TextCell "E17" "  +3.63%  "
Set-TextCell "D18" "42.932.03"
Set-TextCell "E19" "  +7.64%  "
Set-TextCell "D20" "6.19"
Set-TextCell "E20" "  +3.17%  "
Set-TextCell "D21" "0.0₃0906"
Set-TextCell "E21" "  +1.63%  "
Set-TextCell "D22" "67.96"
Set-TextCell "E22" "  +2.17%  "
Set-TextCell "D23" "237.75"
Set-TextCell "E23" "  +1.15%  "
Set-TextCell "D24" "2.21"
Set-TextCell "E24" "  +11.51%  "
Set-TextCell "D25" "2.47"
Set-TextCell "E25" "  +0.35%  "
Set-TextCell "E26" "  -0.14%  "
Set-TextCell "D27" "24.82"
Set-TextCell "E27" "  +3.32%  "
Set-TextCell "B28" "Toncoin"
Set-TextCell "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D28" "2.20"
Set-TextCell "E28" "  +6.55%  "
Set-TextCell "B29" "Monero"
Set-TextCell "C29" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D29" "167.89"
Set-TextCell "E29" "  -0.63%  "
Set-TextCell "D30" "34.18"
Set-TextCell "E30" "  +1.26%  "
Set-TextCell "D31" "9.26"
Set-TextCell "E31" "  +0.55%  "
Set-TextCell "D32" "0.999"
Set-TextCell "E32" "  +0.00%  "
Set-TextCell "E33" "  +2.45%  "
Set-TextCell "E34" "  +1.73%  "
Set-TextCell "E35" "  +3.74%  "
Set-TextCell "D36" "17.02"
Set-TextCell "E36" "  +2.33%  "
Set-TextCell "D37" "0.0693"
Set-TextCell "E37" "  +0.72%  "
Set-TextCell "E38" "  +3.01%  "
Set-TextCell "E39" "  +1.55%  "
Set-TextCell "E40" "  +4.14%  "
Set-TextCell "E41" "  +1.18%  "
Set-TextCell "E42" "  -4.22%  "
Set-TextCell "D43" "2.002.79"
Set-TextCell "E43" "  +1.50%  "
Set-TextCell "E44" "  +3.59%  "
Set-TextCell "E45" "  +7.35%  "
Set-TextCell "D46" "17.58"
Set-TextCell "E46" "  +0.36%  "
Set-TextCell "D47" "2.87"
Set-TextCell "E47" "  +2.82%  "
Set-TextCell "D48" "55.81"
Set-TextCell "E48" "  +6.57%  "
Set-TextCell "D49" "2.526.99"
Set-TextCell "E49" "  +1.13%  "
Set-TextCell "E50" "  +3.17%  "
Set-TextCell "E51" "  +1.33%  "

Write-Output "Applied 86 cell updates"
